# "Elimna EC anteriores y se agregan nuevos, se modifica base de datos"
#
# The "Periodo Mora" / "Valor Mora" database block (rows 16-28, columns E:F)
# is reordered: the oldest-first listing (2003 .. 2103) becomes a
# newest-first listing (2103 .. 2003) - i.e. the 13 data rows are reversed,
# without disturbing the per-row formatting (row 16 keeps the "normal" row
# style, row 28 keeps the special "last row" border style).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 16
$lastRow  = 28
$colPeriodo = 5   # E - Periodo Mora
$colValor   = 6   # F - Valor Mora

$periodos = @()
$valores  = @()
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $periodos += ,$ws.Cells.Item($r, $colPeriodo).Value()
    $valores  += ,$ws.Cells.Item($r, $colValor).Value()
}

$n = $periodos.Count
for ($i = 0; $i -lt $n; $i++) {
    $r = $firstRow + $i
    $ws.Cells.Item($r, $colPeriodo).Value = $periodos[$n - 1 - $i]
    $ws.Cells.Item($r, $colValor).Value   = $valores[$n - 1 - $i]
}
